$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: rename existing items in place
$ws.Range("A3").Value = "Barq's Root Beer"
$ws.Range("A7").Value = "Java Monster Mean Bean"

# Step 2: append new Java Monster flavors at the bottom
$ws.Range("A9").Value = "Java Monster Irish Crème"
$ws.Range("B9").Value = 4
$ws.Range("A10").Value = "Java Monster Café Latte"
$ws.Range("B10").Value = 4
$ws.Range("A11").Value = "Java Monster Loca Moca"
$ws.Range("B11").Value = 4
$ws.Range("A12").Value = "Java Monster Triple Shot French Vanilla"
$ws.Range("B12").Value = 5
$ws.Range("A13").Value = "Java Monster Triple Shot Mocha"
$ws.Range("B13").Value = 5

# Step 3: update costs that changed
$ws.Range("B8").Value = 3

# Step 4: sort data range A2:B13 alphabetically by column A
$rng = $ws.Range("A1:B13")
$rng.Sort($ws.Range("A1"), 1, $null, $null, 1, $null, 1, 1)

$ws.Columns.Item(1).ColumnWidth = 27.6
$ws.Range("B13").Select() | Out-Null
